$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (MuSCs as sending cluster rows removed)
$ws.Rows("8:10").Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.775841999999999
$ws.Range("H2").Value = 8.327525999999999
$ws.Range("I2").Value = 0.0624750527258915
$ws.Range("J2").Value = 0.0624750527258915
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05619066666666667
$ws.Range("N2").Value = 0.168572
$ws.Range("O2").Value = 0.3931387525216601
$ws.Range("P2").Value = 0.39313875252166
$ws.Range("Q2").Value = 0.1559764125413333
$ws.Range("R2").Value = 1.403787712872
$ws.Range("S2").Value = 0.02456136429238192
$ws.Range("T2").Value = 0.02456136429238192

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.775841999999999
$ws.Range("H3").Value = 8.327525999999999
$ws.Range("I3").Value = 0.0624750527258915
$ws.Range("J3").Value = 0.0624750527258915
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.08673766666666667
$ws.Range("N3").Value = 0.260213
$ws.Range("O3").Value = 0.60686124747834
$ws.Range("P3").Value = 0.60686124747834
$ws.Range("Q3").Value = 0.2407700581153333
$ws.Range("R3").Value = 2.166930523038
$ws.Range("S3").Value = 0.03791368843350958
$ws.Range("T3").Value = 0.03791368843350958

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 11.78712033333333
$ws.Range("H4").Value = 35.361361
$ws.Range("I4").Value = 0.2652892219050753
$ws.Range("J4").Value = 0.2652892219050753
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05619066666666667
$ws.Range("N4").Value = 0.168572
$ws.Range("O4").Value = 0.3931387525216601
$ws.Range("P4").Value = 0.39313875252166
$ws.Range("Q4").Value = 0.6623261496102223
$ws.Range("R4").Value = 5.960935346492001
$ws.Range("S4").Value = 0.1042954737572031
$ws.Range("T4").Value = 0.1042954737572031

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.78712033333333
$ws.Range("H5").Value = 35.361361
$ws.Range("I5").Value = 0.2652892219050753
$ws.Range("J5").Value = 0.2652892219050753
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08673766666666667
$ws.Range("N5").Value = 0.260213
$ws.Range("O5").Value = 0.60686124747834
$ws.Range("P5").Value = 0.60686124747834
$ws.Range("Q5").Value = 1.022387314432556
$ws.Range("R5").Value = 9.201485829893002
$ws.Range("S5").Value = 0.1609937481478722
$ws.Range("T5").Value = 0.1609937481478722

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 29.86824466666667
$ws.Range("H6").Value = 89.60473400000001
$ws.Range("I6").Value = 0.6722357253690333
$ws.Range("J6").Value = 0.6722357253690333
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.05619066666666667
$ws.Range("N6").Value = 0.168572
$ws.Range("O6").Value = 0.3931387525216601
$ws.Range("P6").Value = 0.39313875252166
$ws.Range("Q6").Value = 1.678316579983111
$ws.Range("R6").Value = 15.104849219848
$ws.Range("S6").Value = 0.2642819144720751
$ws.Range("T6").Value = 0.264281914472075

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 29.86824466666667
$ws.Range("H7").Value = 89.60473400000001
$ws.Range("I7").Value = 0.6722357253690333
$ws.Range("J7").Value = 0.6722357253690333
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.08673766666666667
$ws.Range("N7").Value = 0.260213
$ws.Range("O7").Value = 0.60686124747834
$ws.Range("P7").Value = 0.60686124747834
$ws.Range("Q7").Value = 2.590701849815778
$ws.Range("R7").Value = 23.316316648342
$ws.Range("S7").Value = 0.4079538108969584
$ws.Range("T7").Value = 0.4079538108969584
